$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "sid21111"
$ws.Range("B3").Value = "sid22222"
$ws.Range("C2").Value = "spw21111"
$ws.Range("C3").Value = "spw22222"

$ws.Range("C3").Select()
